$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.917.95"
$ws.Range("E2").Value = "  -1.28%  "

$ws.Range("D3").Value = "1.638.00"
$ws.Range("E3").Value = "  -0.48%  "

$ws.Range("E4").Value = "  +0.32%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.88%  "

$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("E7").Value = "  +0.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.256"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.60%  "

$ws.Range("E9").Value = "  +0.29%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.62"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.75%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0793"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.02%  "

$ws.Range("D12").Value = "1.864.77"
$ws.Range("E12").Value = "  -0.49%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.26"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.05%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.662.31"
$ws.Range("E14").Value = "  +1.13%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.544"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.37%  "

$ws.Range("D16").Value = "0.0₃0767"
$ws.Range("E16").Value = "  +0.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.96"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.02%  "

$ws.Range("D18").Value = "25.924.61"
$ws.Range("E18").Value = "  -1.20%  "

$ws.Range("E19").Value = "  +0.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.55%  "

$ws.Range("E21").Value = "  -1.37%  "

$ws.Range("E22").Value = "  -1.30%  "

$ws.Range("E23").Value = "  -1.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "143.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.41%  "

$ws.Range("E25").Value = "  +0.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.78"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("E27").Value = "  +1.14%  "

$ws.Range("E28").Value = "  -1.35%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.54%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0502"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.32%  "

$ws.Range("E32").Value = "  -1.53%  "

$ws.Range("E33").Value = "  -0.37%  "

$ws.Range("E34").Value = "  -3.75%  "

$ws.Range("E35").Value = "  +1.26%  "

$ws.Range("E36").Value = "  -1.47%  "

$ws.Range("D37").Value = "1.139.16"
$ws.Range("E37").Value = "  +0.07%  "

$ws.Range("E38").Value = "  -1.81%  "

$ws.Range("E39").Value = "  -2.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0158"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.00%  "

$ws.Range("E41").Value = "  +0.15%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.31"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.88%  "

$ws.Range("E44").Value = "  -0.08%  "

$ws.Range("D45").Value = "1.774.41"
$ws.Range("E45").Value = "  -0.49%  "

$ws.Range("E46").Value = "  +2.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.58%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0533"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.95%  "

$ws.Range("E49").Value = "  -0.78%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.15%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.414"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.73%  "
